$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text is safe to assign directly (Excel will not reinterpret them as numbers) ---
$ws.Range('D2').Value = '62.783.45'
$ws.Range('E2').Value = '  +1.16%  '
$ws.Range('D3').Value = '3.469.02'
$ws.Range('E3').Value = '  +1.14%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('E5').Value = '  +1.15%  '
$ws.Range('E6').Value = '  +0.48%  '
$ws.Range('E7').Value = '  -0.84%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -2.19%  '
$ws.Range('E10').Value = '  +6.26%  '
$ws.Range('E11').Value = '  -0.66%  '
$ws.Range('E12').Value = '  +3.68%  '
$ws.Range('E13').Value = '  -2.80%  '
$ws.Range('D14').Value = '4.022.85'
$ws.Range('E14').Value = '  +1.22%  '
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('D17').Value = '3.477.97'
$ws.Range('E17').Value = '  +0.93%  '
$ws.Range('E18').Value = '  +1.24%  '
$ws.Range('E19').Value = '  -1.77%  '
$ws.Range('D20').Value = '62.742.17'
$ws.Range('E20').Value = '  +1.05%  '
$ws.Range('E21').Value = '  +2.54%  '
$ws.Range('E22').Value = '  -0.89%  '
$ws.Range('E23').Value = '  +1.47%  '
$ws.Range('E24').Value = '  +1.69%  '
$ws.Range('E25').Value = '  +17.02%  '
$ws.Range('E26').Value = '  +1.94%  '
$ws.Range('E27').Value = '  +1.02%  '
$ws.Range('E28').Value = '  +0.70%  '
$ws.Range('E29').Value = '  -2.24%  '
$ws.Range('E30').Value = '  -0.29%  '
$ws.Range('E31').Value = '  -1.63%  '
$ws.Range('E32').Value = '  -2.29%  '
$ws.Range('E33').Value = '  -1.77%  '
$ws.Range('E34').Value = '  -5.04%  '
$ws.Range('E35').Value = '  +0.22%  '
$ws.Range('E36').Value = '  +8.07%  '
$ws.Range('E37').Value = '  -2.23%  '
$ws.Range('B38').Value = 'FirstDigitalUSD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('E38').Value = '  +0.16%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('E39').Value = '  +4.25%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E40').Value = '  +2.72%  '
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('E41').Value = '  -0.78%  '
$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('B43').Value = 'LidoDAOToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('E43').Value = '  -1.31%  '
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('E44').Value = '  +7.18%  '
$ws.Range('E45').Value = '  +3.79%  '
$ws.Range('E46').Value = '  +1.55%  '
$ws.Range('E47').Value = '  +13.60%  '
$ws.Range('E48').Value = '  +28.27%  '
$ws.Range('E49').Value = '  -1.96%  '
$ws.Range('E50').Value = '  -0.15%  '
$ws.Range('E51').Value = '  +0.68%  '

# --- D-column values that look like plain numbers (e.g. "1.00", "9.60"): force Text format first so
#     Excel keeps the exact original formatting (trailing zeros, etc.) instead of coercing to a number,
#     then clear the temporary number format so the cell style matches the untouched cells again. ---
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '414.07'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.62'
$ws.Range('D6').ClearFormats()
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.628'
$ws.Range('D7').ClearFormats()
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').ClearFormats()
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.151'
$ws.Range('D10').ClearFormats()
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.61'
$ws.Range('D11').ClearFormats()
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '9.60'
$ws.Range('D12').ClearFormats()
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000222'
$ws.Range('D13').ClearFormats()
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.140'
$ws.Range('D15').ClearFormats()
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.52'
$ws.Range('D16').ClearFormats()
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.67'
$ws.Range('D18').ClearFormats()
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '465.53'
$ws.Range('D21').ClearFormats()
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '90.92'
$ws.Range('D22').ClearFormats()
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.27'
$ws.Range('D23').ClearFormats()
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.56'
$ws.Range('D25').ClearFormats()
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.31'
$ws.Range('D26').ClearFormats()
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '33.41'
$ws.Range('D27').ClearFormats()
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.56'
$ws.Range('D29').ClearFormats()
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '12.00'
$ws.Range('D30').ClearFormats()
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.65'
$ws.Range('D31').ClearFormats()
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '40.78'
$ws.Range('D34').ClearFormats()
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '58.54'
$ws.Range('D36').ClearFormats()
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0492'
$ws.Range('D37').ClearFormats()
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('D38').ClearFormats()
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.07'
$ws.Range('D39').ClearFormats()
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '146.95'
$ws.Range('D40').ClearFormats()
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.134'
$ws.Range('D41').ClearFormats()
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.321'
$ws.Range('D42').ClearFormats()
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.33'
$ws.Range('D43').ClearFormats()
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.70'
$ws.Range('D44').ClearFormats()
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.40'
$ws.Range('D47').ClearFormats()
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '16.37'
$ws.Range('D49').ClearFormats()
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.34'
$ws.Range('D50').ClearFormats()
